function Update-EventSheet($ws) {
    # 1. Direct F-column numeric updates (rows unaffected by later insert)
    $ws.Range("F4").Value = 48
    $ws.Range("F5").Value = 248
    $ws.Range("F7").Value = 126
    $ws.Range("F14").Value = 363
    $ws.Range("F16").Value = 467
    $ws.Range("F17").Value = 393
    $ws.Range("F18").Value = 134
    $ws.Range("F19").Value = 61
    $ws.Range("F21").Value = 38
    $ws.Range("F22").Value = 985
    $ws.Range("F23").Value = 2757
    $ws.Range("F26").Value = 526

    # 2. Insert a new row at 27, shifting old rows 27-34 down to 28-35
    $ws.Rows.Item(27).Insert()

    # Fix up formatting of new A27 (copy format only from A28, which now carries old row27's style)
    $ws.Range("A28").Copy()
    $ws.Range("A27").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $ws.Range("A27").Value = 26

    # 3. Populate new row 27 with the 吉安 event
    $ws.Range("B27").NumberFormat = "@"
    $ws.Range("B27").Value = "2024-08-03"
    $ws.Range("C27").Value = "吉安·COMIC LIFE周年庆典"
    $ws.Range("D27").Value = "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心"
    $ws.Range("E27").Value = "2024.08.03 09:30-08.03 18:00"
    $ws.Range("F27").Value = 14
    $ws.Range("G27").Value = 9.9
    $ws.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=87164"
    $ws.Range("I27").Value = "//i1.hdslb.com/bfs/openplatform/202406/NWD9iQ9h1717598526259.jpeg"

    # 4. Update F-values for the shifted rows that also changed content
    $ws.Range("F28").Value = 969
    $ws.Range("F31").Value = 258
    $ws.Range("F32").Value = 385
}

$wb = $excel.ActiveWorkbook
Update-EventSheet($wb.Worksheets.Item("展览"))
Update-EventSheet($wb.Worksheets.Item("全部类型"))
Write-Host "done"
